$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 456.64706
$ws.Range("I19").Value = 387.6
$ws.Range("J19").Value = 511.1579
$ws.Range("K19").Value = 387.6
$ws.Range("L19").Value = 511.1579
$ws.Range("M19").Value = -212.6
$ws.Range("N19").Value = -861.1578999999999
$ws.Range("H21").Value = 14000
$ws.Range("I21").Value = 14000
$ws.Range("K21").Value = 14000
$ws.Range("M21").Value = -13532
$ws.Range("H23").Value = 14000
$ws.Range("I23").Value = 14000
$ws.Range("K23").Value = 14000
$ws.Range("M23").Value = -13766
$ws.Range("H42").Value = 1251.3
$ws.Range("I42").Value = 1943.8334
$ws.Range("J42").Value = 212.5
$ws.Range("K42").Value = 5831.5002
$ws.Range("L42").Value = 637.5
$ws.Range("M42").Value = -5601.5002
$ws.Range("N42").Value = -1097.5
$ws.Range("H43").Value = 1087.875
$ws.Range("I43").Value = 1083.5
$ws.Range("J43").Value = 1101
$ws.Range("K43").Value = 1083.5
$ws.Range("L43").Value = 1101
$ws.Range("M43").Value = -1014.5
$ws.Range("N43").Value = -1239
$ws.Range("H62").Value = 3102654.5
$ws.Range("I62").Value = 4128125.5
$ws.Range("J62").Value = 26241.889
$ws.Range("K62").Value = 4128125.5
$ws.Range("L62").Value = 26241.889
$ws.Range("M62").Value = -4127501.5
$ws.Range("N62").Value = -27489.889
$ws.Range("H65").Value = 3102654.5
$ws.Range("I65").Value = 4128125.5
$ws.Range("J65").Value = 26241.889
$ws.Range("K65").Value = 20640627.5
$ws.Range("L65").Value = 131209.445
$ws.Range("M65").Value = -20637507.5
$ws.Range("N65").Value = -137449.445
$ws.Range("H87").Value = 35660.5
$ws.Range("J87").Value = 59000
$ws.Range("L87").Value = 59000
$ws.Range("N87").Value = -61496
$ws.Range("H90").Value = 35660.5
$ws.Range("J90").Value = 59000
$ws.Range("L90").Value = 177000
$ws.Range("N90").Value = -189480
$ws.Range("H107").Value = 505457.03
$ws.Range("I107").Value = 555927.9399999999
$ws.Range("J107").Value = 748
$ws.Range("K107").Value = 555927.9399999999
$ws.Range("L107").Value = 748
$ws.Range("M107").Value = -554007.9399999999
$ws.Range("N107").Value = -4588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 335633.66
$ws.Range("I5").Value = 335633.66
$ws.Range("K5").Value = 335633.66
$ws.Range("M5").Value = -335521.66
$ws.Range("H41").Value = 928
$ws.Range("I41").Value = 928
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 928
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -514
$ws.Range("N41").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 335633.66
$ws.Range("I4").Value = 335633.66
$ws.Range("K4").Value = 335633.66
$ws.Range("M4").Value = -335518.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 21601.666
$ws.Range("I2").Value = 29900
$ws.Range("J2").Value = 5005
$ws.Range("K2").Value = 29900
$ws.Range("L2").Value = 5005
$ws.Range("M2").Value = -29787
$ws.Range("N2").Value = -5231
$ws.Range("H4").Value = 1000000000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000000000
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -1000000224
$ws.Range("H22").Value = 994.875
$ws.Range("I22").Value = 661.5
$ws.Range("J22").Value = 1995
$ws.Range("K22").Value = 661.5
$ws.Range("L22").Value = 1995
$ws.Range("M22").Value = -311.5
$ws.Range("N22").Value = -2695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4006.6667
$ws.Range("I69").Value = 600
$ws.Range("J69").Value = 4858.3335
$ws.Range("K69").Value = 1800
$ws.Range("L69").Value = 14575.0005
$ws.Range("M69").Value = -989
$ws.Range("N69").Value = -16197.0005
$ws.Range("H72").Value = 4006.6667
$ws.Range("I72").Value = 600
$ws.Range("J72").Value = 4858.3335
$ws.Range("K72").Value = 5400
$ws.Range("L72").Value = 43725.0015
$ws.Range("M72").Value = -1344
$ws.Range("N72").Value = -51837.0015
$ws.Range("H104").Value = 3676.3333
$ws.Range("J104").Value = 3676.3333
$ws.Range("L104").Value = 11028.9999
$ws.Range("N104").Value = -16270.9999
$ws.Range("H119").Value = 1046.091
$ws.Range("I119").Value = 656.3333
$ws.Range("K119").Value = 1968.9999
$ws.Range("M119").Value = 2869.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("N21").Value = $null
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("N30").Value = $null
$ws.Range("H107").Value = 1295.3334
$ws.Range("I107").Value = 1004.4
$ws.Range("K107").Value = 1004.4
$ws.Range("M107").Value = 915.6
$ws.Range("H123").Value = 11031.143
$ws.Range("J123").Value = 11031.143
$ws.Range("L123").Value = 11031.143
$ws.Range("N123").Value = -15931.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3199.75
$ws.Range("I9").Value = 3199.75
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3199.75
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -2975.75
$ws.Range("N9").Value = $null
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H35").Value = 24750
$ws.Range("I35").Value = 45000
$ws.Range("K35").Value = 45000
$ws.Range("M35").Value = -44664
$ws.Range("H46").Value = 715.35297
$ws.Range("I46").Value = 560.9091
$ws.Range("J46").Value = 998.5
$ws.Range("K46").Value = 560.9091
$ws.Range("L46").Value = 998.5
$ws.Range("M46").Value = -372.9091
$ws.Range("N46").Value = -1374.5
$ws.Range("H61").Value = 735.3043
$ws.Range("I61").Value = 461.93332
$ws.Range("J61").Value = 1247.875
$ws.Range("K61").Value = 461.93332
$ws.Range("L61").Value = 1247.875
$ws.Range("M61").Value = -259.93332
$ws.Range("N61").Value = -1651.875
$ws.Range("H113").Value = 735.3043
$ws.Range("I113").Value = 461.93332
$ws.Range("J113").Value = 1247.875
$ws.Range("K113").Value = 461.93332
$ws.Range("L113").Value = 1247.875
$ws.Range("M113").Value = 1708.06668
$ws.Range("N113").Value = -5587.875
$ws.Range("H127").Value = 48500
$ws.Range("J127").Value = 48500
$ws.Range("L127").Value = 48500
$ws.Range("N127").Value = -58420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 7749.25
$ws.Range("I17").Value = 7749.25
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 7749.25
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -7577.25
$ws.Range("N17").Value = $null
$ws.Range("H26").Value = 9666.666999999999
$ws.Range("J26").Value = 9666.666999999999
$ws.Range("L26").Value = 9666.666999999999
$ws.Range("N26").Value = -10252.667
$ws.Range("H122").Value = 1379.2
$ws.Range("I122").Value = 1379.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1687.6
$ws.Range("N122").Value = $null
